$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a numeric-looking value as TEXT (shared string) into a cell,
# while preserving that cell's existing style (so no new cellXfs entries are
# created). Trick: build the text via a formula in a scratch cell, copy it,
# then PasteSpecial only the *values* into the destination - this keeps the
# destination cell's current formatting/style untouched.
# ---------------------------------------------------------------------------
function Set-TextValue($addr, $text) {
    $ws.Range("ZZ1").Formula = '="' + $text + '"'
    $ws.Range("ZZ1").Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $ws.Range("ZZ1").Value = ""
}

# ---------------------------------------------------------------------------
# Row 3: header row (G1 / G2 run labels)
# ---------------------------------------------------------------------------
$ws.Range("F3").Value = "G2 #1"
$ws.Range("G3").Value = "G2 #2"
$ws.Range("H3").Value = "G2 #3"
$ws.Range("I3").Value = "G2 #4"

# ---------------------------------------------------------------------------
# NODE EMBEDDING block (rows 4-9), columns G,H,I are new
# ---------------------------------------------------------------------------
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 5
$ws.Range("I4").Value = 10

$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 6
$ws.Range("I5").Value = 6

$ws.Range("G6").Value = 7
$ws.Range("H6").Value = 7
$ws.Range("I6").Value = 5

$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0.75
$ws.Range("I7").Value = 0.5

$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1

$ws.Range("G9").Value = 128
$ws.Range("H9").Value = 64
$ws.Range("I9").Value = 64

# ---------------------------------------------------------------------------
# LINK PREDICTION GNN block (rows 14-20), columns G,H,I are new
# ---------------------------------------------------------------------------
$ws.Range("G14").Value = 256
$ws.Range("H14").Value = 256
$ws.Range("I14").Value = 64

$ws.Range("G15").Value = 256
$ws.Range("H15").Value = 64
$ws.Range("I15").Value = 64

$ws.Range("G16").Value = 200
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 150

# row 17 (lr) - these text values look numeric, must remain text
Set-TextValue "G17" "0.017630062959213853"
Set-TextValue "H17" "0.013477408995651594"
Set-TextValue "I17" "0.00040232263806239127"

# row 18 (aggregation function)
$ws.Range("G18").Value = "mean"
$ws.Range("H18").Value = "mean"
$ws.Range("I18").Value = "mean"

# row 19 (dropout)
$ws.Range("G19").Value = 0.1
$ws.Range("H19").Value = 0.2
$ws.Range("I19").Value = 0.1

# row 20 (layers)
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 2
$ws.Range("I20").Value = 2

# ---------------------------------------------------------------------------
# SCORES block (rows 23-25)
# ---------------------------------------------------------------------------
# row 23 (train) - E23 value updated, F-I are new
$ws.Range("E23").Value = 0.99029999999999996
$ws.Range("F23").Value = 0.98750000000000004
$ws.Range("G23").Value = 0.98419999999999996
$ws.Range("H23").Value = 0.98250000000000004
$ws.Range("I23").Value = 0.95740000000000003

# row 24 (validation) - E24 becomes text "-", F-I are new
Set-TextValue "E24" "-"
$ws.Range("F24").Value = 0.96150000000000002
$ws.Range("G24").Value = 0.97699999999999998
$ws.Range("H24").Value = 0.97440000000000004
$ws.Range("I24").Value = 0.94650000000000001

# row 25 (test) - E25 value updated, F-I are new
$ws.Range("E25").Value = 0.98250000000000004
$ws.Range("F25").Value = 0.96419999999999995
$ws.Range("G25").Value = 0.9758
$ws.Range("H25").Value = 0.97840000000000005
$ws.Range("I25").Value = 0.9476

# ---------------------------------------------------------------------------
# Apply consistent formatting (style) to every newly-populated cell in
# columns G, H, I by copying the format from an existing centre-aligned
# data cell (E19, style index 1) - this re-uses the existing style instead
# of creating brand-new cellXfs records.
# ---------------------------------------------------------------------------
$ws.Range("E19").Copy()
$ws.Range("G3:I9").PasteSpecial(-4122)
$ws.Range("G14:I16").PasteSpecial(-4122)
$ws.Range("G17:H17").PasteSpecial(-4122)
$ws.Range("G18:I20").PasteSpecial(-4122)
$ws.Range("F23:I23").PasteSpecial(-4122)
$ws.Range("E24").PasteSpecial(-4122)
$ws.Range("F24:I24").PasteSpecial(-4122)
$ws.Range("F25:I25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Column widths for the new columns G, H, I
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 25.59
$ws.Columns.Item(8).ColumnWidth = 20.75
$ws.Columns.Item(9).ColumnWidth = 21.1

# ---------------------------------------------------------------------------
# Selection, matching the saved view state in the target workbook
# ---------------------------------------------------------------------------
$ws.Range("I25").Select()
